$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("g10.2")

$values = @{
    "B2"  = -3.383341250834726
    "C2"  = 2.155784801988059
    "D2"  = 7.433110380454799
    "B3"  = 0.2439808399159471
    "C3"  = -0.1740913406866107
    "D3"  = -1.008724857903942
    "B4"  = 5.606362531019382
    "C4"  = 0.6579095658640677
    "D4"  = 6.266452669837808
    "B5"  = 6.187463753814892
    "C5"  = -6.903561235429456
    "D5"  = 9.754254968344501
    "B6"  = -1.451057335467054
    "C6"  = -6.177742879580272
    "D6"  = 8.002056008584525
    "B7"  = -0.2896224593916363
    "C7"  = -4.969367067900111
    "D7"  = 2.803087220402856
    "B8"  = -0.9081980491532082
    "C8"  = -4.256032038651048
    "D8"  = 0.1156072649098894
    "B9"  = 4.642033376017518
    "C9"  = -1.180271120538434
    "D9"  = 11.24620493149768
    "B10" = -10.64870586592827
    "C10" = -5.575599503798379
    "D10" = -5.924122471118087
    "B11" = -6.412425571460922
    "C11" = 8.995566678924227
    "D11" = -7.347264789576724
    "B12" = 0.08943577057456409
    "C12" = 7.480578958599393
    "D12" = -10.81946663818691
    "B13" = -2.468191820892673
    "C13" = 2.943460101589301
    "D13" = -3.206398799120913
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
